$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing translations
# Row 104: "Configuration" row - change Hebrew translation from "תצורה" to "הגדרות"
$ws.Range("C104").Value = "הגדרות"

# Row 102: "Orders" row - correct Hebrew translation from "משלוחים" (Deliveries) to "הזמנות" (Orders)
$ws.Range("C102").Value = "הזמנות"

# Add a new row 114 for "Attributes"
$ws.Range("B114").Value = "Attributes"
$ws.Range("C114").Value = "תוספות למוצרים"

# Update the view so the new row is visible, matching author's scroll position
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F119").Select()
